$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.350.32"
$ws.Range("D3").Value = "3.139.52"
$ws.Range("E3").Value = "  -4.91%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'516.89"
$ws.Range("E5").Value = "  -7.13%  "
$ws.Range("D6").Value = "'133.18"
$ws.Range("E6").Value = "  -7.08%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "3.141.56"
$ws.Range("E8").Value = "  -4.96%  "
$ws.Range("E9").Value = "  -6.37%  "
$ws.Range("D10").Value = "'7.20"
$ws.Range("E10").Value = "  -8.35%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  -9.56%  "
$ws.Range("D12").Value = "'0.381"
$ws.Range("E12").Value = "  -6.32%  "
$ws.Range("D13").Value = "3.669.18"
$ws.Range("E13").Value = "  -4.82%  "
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "'25.15"
$ws.Range("E15").Value = "  -7.05%  "
$ws.Range("D16").Value = "3.136.05"
$ws.Range("E16").Value = "  -4.91%  "
$ws.Range("D17").Value = "57.300.76"
$ws.Range("E17").Value = "  -4.73%  "
$ws.Range("E18").Value = "  -10.49%  "
$ws.Range("D19").Value = "'5.71"
$ws.Range("E19").Value = "  -6.66%  "
$ws.Range("E20").Value = "  -9.92%  "
$ws.Range("D21").Value = "'7.92"
$ws.Range("E21").Value = "  -7.46%  "
$ws.Range("D22").Value = "'342.16"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'68.33"
$ws.Range("E24").Value = "  -7.17%  "
$ws.Range("D25").Value = "'0.501"
$ws.Range("E25").Value = "  -7.94%  "
$ws.Range("D26").Value = "3.260.80"
$ws.Range("E26").Value = "  -5.15%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.163"
$ws.Range("E28").Value = "  -5.91%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0927"
$ws.Range("E29").Value = "  -9.88%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -7.15%  "
$ws.Range("D32").Value = "'1.83"
$ws.Range("E32").Value = "  -9.19%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'21.48"
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'6.81"
$ws.Range("E34").Value = "  -10.79%  "
$ws.Range("E35").Value = "  -5.30%  "
$ws.Range("D36").Value = "'4.81"
$ws.Range("E36").Value = "  -7.32%  "
$ws.Range("D37").Value = "'157.39"
$ws.Range("E37").Value = "  -5.50%  "
$ws.Range("D38").Value = "'6.15"
$ws.Range("E38").Value = "  -8.58%  "
$ws.Range("E39").Value = "  -9.77%  "
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("D41").Value = "3.166.07"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").Value = "'0.0680"
$ws.Range("E42").Value = "  -8.24%  "
$ws.Range("D43").Value = "'40.22"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("D44").Value = "'0.690"
$ws.Range("E44").Value = "  -8.03%  "
$ws.Range("D45").Value = "'1.06"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'3.84"
$ws.Range("E47").Value = "  -8.51%  "
$ws.Range("E48").Value = "  -9.21%  "
$ws.Range("D49").Value = "2.228.07"
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("D50").Value = "'6.10"
$ws.Range("E50").Value = "  -6.49%  "
$ws.Range("D51").Value = "'19.86"
$ws.Range("E51").Value = "  -6.10%  "
